{"js": "const body = context.document.body;\n\nconst pairs = [\n  [\"2025-09-10 Wednesday\", \"2025-09-11 Thursday\"],\n  [\"553\u00d75=2765\", \"306\u00d75=1530\"],\n  [\"523\u00d73=1569\", \"383\u00d74=1532\"],\n  [\"483\u00d73=1449\", \"676\u00d77=4732\"],\n  [\"720\u00d73=2160\", \"810\u00d79=7290\"],\n  [\"208\u00d76=1248\", \"880\u00d79=7920\"],\n  [\"555\u00d73=1665\", \"607\u00d73=1821\"],\n  [\"427\u00d74=1708\", \"356\u00d75=1780\"],\n  [\"843\u00d76=5058\", \"522\u00d77=3654\"],\n  [\"633\u00d73=1899\", \"500\u00d76=3000\"],\n  [\"369\u00d79=3321\", \"667\u00d78=5336\"],\n  [\"649\u00d72=1298\", \"394\u00d72=788\"],\n  [\"551\u00d78=4408\", \"644\u00d75=3220\"],\n  [\"938\u00d74=3752\", \"189\u00d78=1512\"],\n  [\"156\u00d74=624\", \"920\u00d73=2760\"],\n  [\"731\u00d74=2924\", \"129\u00d73=387\"],\n  [\"624\u00d73=1872\", \"780\u00d74=3120\"],\n  [\"167\u00d75=835\", \"681\u00d74=2724\"],\n  [\"529\u00d75=2645\", \"894\u00d78=7152\"],\n  [\"575\u00d72=1150\", \"686\u00d73=2058\"],\n  [\"165\u00d78=1320\", \"509\u00d73=1527\"],\n  [\"736\u00d77=5152\", \"879\u00d79=7911\"],\n  [\"703\u00d76=4218\", \"457\u00d79=4113\"],\n  [\"440\u00d75=2200\", \"196\u00d72=392\"],\n  [\"191\u00d76=1146\", \"220\u00d73=660\"],\n  [\"528\u00d72=1056\", \"723\u00d77=5061\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"2025-09-10 Wednesday\", \"2025-09-11 Thursday\")\n  ,@(\"553\u00d75=2765\", \"306\u00d75=1530\")\n  ,@(\"523\u00d73=1569\", \"383\u00d74=1532\")\n  ,@(\"483\u00d73=1449\", \"676\u00d77=4732\")\n  ,@(\"720\u00d73=2160\", \"810\u00d79=7290\")\n  ,@(\"208\u00d76=1248\", \"880\u00d79=7920\")\n  ,@(\"555\u00d73=1665\", \"607\u00d73=1821\")\n  ,@(\"427\u00d74=1708\", \"356\u00d75=1780\")\n  ,@(\"843\u00d76=5058\", \"522\u00d77=3654\")\n  ,@(\"633\u00d73=1899\", \"500\u00d76=3000\")\n  ,@(\"369\u00d79=3321\", \"667\u00d78=5336\")\n  ,@(\"649\u00d72=1298\", \"394\u00d72=788\")\n  ,@(\"551\u00d78=4408\", \"644\u00d75=3220\")\n  ,@(\"938\u00d74=3752\", \"189\u00d78=1512\")\n  ,@(\"156\u00d74=624\", \"920\u00d73=2760\")\n  ,@(\"731\u00d74=2924\", \"129\u00d73=387\")\n  ,@(\"624\u00d73=1872\", \"780\u00d74=3120\")\n  ,@(\"167\u00d75=835\", \"681\u00d74=2724\")\n  ,@(\"529\u00d75=2645\", \"894\u00d78=7152\")\n  ,@(\"575\u00d72=1150\", \"686\u00d73=2058\")\n  ,@(\"165\u00d78=1320\", \"509\u00d73=1527\")\n  ,@(\"736\u00d77=5152\", \"879\u00d79=7911\")\n  ,@(\"703\u00d76=4218\", \"457\u00d79=4113\")\n  ,@(\"440\u00d75=2200\", \"196\u00d72=392\")\n  ,@(\"191\u00d76=1146\", \"220\u00d73=660\")\n  ,@(\"528\u00d72=1056\", \"723\u00d77=5061\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}"}
